$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "Tổng xu" (H8) becomes "Tổng doanh thu", and the old
# "Tổng doanh thu" (I8) becomes "Tổng doanh thu xu" (new point/coin revenue column).
$ws.Range("H8").Value = "Tổng doanh thu"
$ws.Range("I8").Value = "Tổng doanh thu xu"

$ws.Range("D14").Select()
